# Apply weekly re-shuffle of the Fecha/Volumen/Precio/Origen columns
# (columns D, J, K, L, M, O, P) across data rows 2-14.
#
# The underlying data rows stay in place for the "static" columns
# (A, B, C, E, F, G, H, I, N, Q, R are identical for every row already)
# but the per-record values for Fecha (D), Volumen (J), Precio minimo/
# maximo/promedio (K/L/M), Origen (O) and Precio $/Kg (P) get
# redistributed among the rows, following this row->row mapping
# (source row number -> destination row number):
#
#   2 -> 3, 3 -> 14, 4 -> 13, 5 -> 11, 6 -> 2, 7 -> 9, 8 -> 12,
#   9 -> 7, 10 -> 8, 11 -> 5, 12 -> 4, 13 -> 6, 14 -> 10

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$mapping = @{
    2  = 3
    3  = 14
    4  = 13
    5  = 11
    6  = 2
    7  = 9
    8  = 12
    9  = 7
    10 = 8
    11 = 5
    12 = 4
    13 = 6
    14 = 10
}

$cols = @("D", "J", "K", "L", "M", "O", "P")

# Snapshot all the current values for the affected columns/rows before
# writing anything, since this is a permutation (not simple pairwise
# swaps) and rows must not read already-overwritten data.
$snapshot = @{}
foreach ($row in $mapping.Keys) {
    $rowValues = @{}
    foreach ($col in $cols) {
        $rowValues[$col] = $ws.Range("$col$row").Value2
    }
    $snapshot[$row] = $rowValues
}

foreach ($srcRow in $mapping.Keys) {
    $destRow = $mapping[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value2 = $snapshot[$srcRow][$col]
    }
}
